# Update the crypto price/volume snapshot (columns D and E) for rows 2-51.
# Values that look like plain numbers are prefixed with a leading apostrophe
# so Excel stores them as text (matching the original inlineStr cells)
# instead of silently converting them to numeric values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.682.02"
$ws.Range("E2").Value = "  -0.91%  "
$ws.Range("D3").Value = "3.785.72"
$ws.Range("E3").Value = "  +1.20%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'595.53"
$ws.Range("E5").Value = "  +0.44%  "
$ws.Range("D6").Value = "'167.10"
$ws.Range("E6").Value = "  +0.71%  "
$ws.Range("D7").Value = "3.771.50"
$ws.Range("E7").Value = "  +0.93%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  +0.34%  "
$ws.Range("E10").Value = "  +0.03%  "
$ws.Range("D11").Value = "'6.29"
$ws.Range("E11").Value = "  -2.10%  "
$ws.Range("E12").Value = "  +0.28%  "
$ws.Range("E13").Value = "  -2.71%  "
$ws.Range("D14").Value = "'35.96"
$ws.Range("E14").Value = "  -0.19%  "
$ws.Range("D15").Value = "4.421.95"
$ws.Range("E15").Value = "  +1.17%  "
$ws.Range("D16").Value = "3.775.63"
$ws.Range("E16").Value = "  +0.73%  "
$ws.Range("D17").Value = "'18.60"
$ws.Range("E17").Value = "  +4.60%  "
$ws.Range("D18").Value = "67.639.62"
$ws.Range("E18").Value = "  -0.95%  "
$ws.Range("E19").Value = "  +0.99%  "
$ws.Range("E20").Value = "  -0.11%  "
$ws.Range("E21").Value = "  -5.57%  "
$ws.Range("E22").Value = "  -1.08%  "
$ws.Range("E23").Value = "  +0.44%  "
$ws.Range("E24").Value = "  +4.98%  "
$ws.Range("D25").Value = "'83.40"
$ws.Range("E25").Value = "  -0.57%  "
$ws.Range("D26").Value = "'11.97"
$ws.Range("E26").Value = "  +1.01%  "
$ws.Range("D27").Value = "'2.12"
$ws.Range("E27").Value = "  -2.40%  "
$ws.Range("E28").Value = "  +0.15%  "
$ws.Range("E29").Value = "  -0.27%  "
$ws.Range("D30").Value = "3.931.49"
$ws.Range("E30").Value = "  +1.07%  "
$ws.Range("E31").Value = "  +0.53%  "
$ws.Range("E32").Value = "  +3.89%  "
$ws.Range("E33").Value = "  -1.14%  "
$ws.Range("D34").Value = "'29.63"
$ws.Range("E34").Value = "  -0.45%  "
$ws.Range("D35").Value = "'0.999"
$ws.Range("E35").Value = "  +0.05%  "
$ws.Range("E36").Value = "  -0.45%  "
$ws.Range("E37").Value = "  -0.56%  "
$ws.Range("E38").Value = "  -1.67%  "
$ws.Range("E39").Value = "  -0.36%  "
$ws.Range("D40").Value = "'0.995"
$ws.Range("E40").Value = "  +0.17%  "
$ws.Range("E41").Value = "  +0.01%  "
$ws.Range("D42").Value = "'0.999"
$ws.Range("E42").Value = "  -0.07%  "
$ws.Range("E43").Value = "  -0.01%  "
$ws.Range("D44").Value = "'45.55"
$ws.Range("E44").Value = "  +5.82%  "
$ws.Range("D45").Value = "'48.15"
$ws.Range("E45").Value = "  +3.10%  "
$ws.Range("E46").Value = "  -0.51%  "
$ws.Range("D47").Value = "'150.24"
$ws.Range("E47").Value = "  +4.13%  "
$ws.Range("D48").Value = "'8.31"
$ws.Range("E48").Value = "  -1.44%  "
$ws.Range("D49").Value = "'392.69"
$ws.Range("D50").Value = "'26.62"
$ws.Range("E50").Value = "  +6.40%  "
$ws.Range("D51").Value = "'1.81"
$ws.Range("E51").Value = "  -5.12%  "
